# Generate Report for Handoff
# Adds two new localization entries (c4713302-... and e55b63a1-...) to the
# Overview sheet as well as the zh-cn and de-de detail sheets, growing each
# table/range from 3 rows (A1:x3) to 5 rows (A1:x5).
#
# Note: values such as "True"/"False" and the blank placeholder need a
# leading apostrophe so the engine stores them as text (matching the
# original data, which uses shared strings rather than native booleans),
# instead of being auto-coerced to Boolean values.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"
$blank = "'"

# ---------------------------------------------------------------------
# Sheet 1: Overview (table "Overview", columns A:G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4 - c4713302-cd18-4021-9747-9a8f13518cb8.md
$wsOverview.Cells.Item(4, 1).Value = "c4713302-cd18-4021-9747-9a8f13518cb8.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c4713302cd184021974400000000000000000000/e2e/c4713302-cd18-4021-9747-9a8f13518cb8.md", "", "", "e2e\c4713302-cd18-4021-9747-9a8f13518cb8.md")
$wsOverview.Cells.Item(4, 3).Value = ".md"
$wsOverview.Cells.Item(4, 4).Value = $blank
$wsOverview.Cells.Item(4, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 7).NumberFormat = $dateFmt
$wsOverview.Cells.Item(4, 7).Value = "2016-08-12 20:49:09"

# Row 5 - e55b63a1-52a7-4daf-a112-431e7e259a2c.md
$wsOverview.Cells.Item(5, 1).Value = "e55b63a1-52a7-4daf-a112-431e7e259a2c.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e55b63a152a74daf00000000000000000000000/e2e/e55b63a1-52a7-4daf-a112-431e7e259a2c.md", "", "", "e2e\e55b63a1-52a7-4daf-a112-431e7e259a2c.md")
$wsOverview.Cells.Item(5, 3).Value = ".md"
$wsOverview.Cells.Item(5, 4).Value = $blank
$wsOverview.Cells.Item(5, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(5, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(5, 7).NumberFormat = $dateFmt
$wsOverview.Cells.Item(5, 7).Value = "2016-08-12 20:49:09"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn (table "zh-cn" / displayName zh_cn, columns A:P)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

# Row 4 - c4713302-cd18-4021-9747-9a8f13518cb8.md
$wsZhCn.Cells.Item(4, 1).Value = "c4713302-cd18-4021-9747-9a8f13518cb8.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c4713302cd184021974400000000000000000000/e2e/c4713302-cd18-4021-9747-9a8f13518cb8.md", "", "", "c4713302-cd18-4021-9747-9a8f13518cb8.md")
$wsZhCn.Cells.Item(4, 2).Value = ".md"
$wsZhCn.Cells.Item(4, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(4, 4).Value = "e2e"
$wsZhCn.Cells.Item(4, 5).Value = "ht"
$wsZhCn.Cells.Item(4, 6).Value = "'False"
$wsZhCn.Cells.Item(4, 7).Value = "c4713302-cd18-4021-9747-9a8f13518cb8.452252d09bd25f7a44047af49243f6dc5f980c53.zh-cn.xlf"
$wsZhCn.Cells.Item(4, 8).NumberFormat = $dateFmt
$wsZhCn.Cells.Item(4, 8).Value = "2016-08-12 20:48:56"
$wsZhCn.Cells.Item(4, 9).Value = $blank
$wsZhCn.Cells.Item(4, 10).Value = $blank
$wsZhCn.Cells.Item(4, 11).NumberFormat = $dateFmt
$wsZhCn.Cells.Item(4, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(4, 12).Value = $blank
$wsZhCn.Cells.Item(4, 13).Value = "'True"
$wsZhCn.Cells.Item(4, 14).Value = $blank
$wsZhCn.Cells.Item(4, 15).Value = "'False"
$wsZhCn.Cells.Item(4, 16).Value = $blank

# Row 5 - e55b63a1-52a7-4daf-a112-431e7e259a2c.md
$wsZhCn.Cells.Item(5, 1).Value = "e55b63a1-52a7-4daf-a112-431e7e259a2c.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e55b63a152a74daf00000000000000000000000/e2e/e55b63a1-52a7-4daf-a112-431e7e259a2c.md", "", "", "e55b63a1-52a7-4daf-a112-431e7e259a2c.md")
$wsZhCn.Cells.Item(5, 2).Value = ".md"
$wsZhCn.Cells.Item(5, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(5, 4).Value = "e2e"
$wsZhCn.Cells.Item(5, 5).Value = "ht"
$wsZhCn.Cells.Item(5, 6).Value = "'False"
$wsZhCn.Cells.Item(5, 7).Value = "e55b63a1-52a7-4daf-a112-431e7e259a2c.f3fd377ffe290dc282c6d067ebb59dc8395379ee.zh-cn.xlf"
$wsZhCn.Cells.Item(5, 8).NumberFormat = $dateFmt
$wsZhCn.Cells.Item(5, 8).Value = "2016-08-12 20:48:56"
$wsZhCn.Cells.Item(5, 9).Value = $blank
$wsZhCn.Cells.Item(5, 10).Value = $blank
$wsZhCn.Cells.Item(5, 11).NumberFormat = $dateFmt
$wsZhCn.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(5, 12).Value = $blank
$wsZhCn.Cells.Item(5, 13).Value = "'True"
$wsZhCn.Cells.Item(5, 14).Value = $blank
$wsZhCn.Cells.Item(5, 15).Value = "'False"
$wsZhCn.Cells.Item(5, 16).Value = $blank

# ---------------------------------------------------------------------
# Sheet 3: de-de (table "de-de" / displayName de_de, columns A:P)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

# Row 4 - c4713302-cd18-4021-9747-9a8f13518cb8.md
$wsDeDe.Cells.Item(4, 1).Value = "c4713302-cd18-4021-9747-9a8f13518cb8.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/c4713302cd184021974400000000000000000000/e2e/c4713302-cd18-4021-9747-9a8f13518cb8.md", "", "", "c4713302-cd18-4021-9747-9a8f13518cb8.md")
$wsDeDe.Cells.Item(4, 2).Value = ".md"
$wsDeDe.Cells.Item(4, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(4, 4).Value = "e2e"
$wsDeDe.Cells.Item(4, 5).Value = "ht"
$wsDeDe.Cells.Item(4, 6).Value = "'False"
$wsDeDe.Cells.Item(4, 7).Value = "c4713302-cd18-4021-9747-9a8f13518cb8.452252d09bd25f7a44047af49243f6dc5f980c53.de-de.xlf"
$wsDeDe.Cells.Item(4, 8).NumberFormat = $dateFmt
$wsDeDe.Cells.Item(4, 8).Value = "2016-08-12 20:49:09"
$wsDeDe.Cells.Item(4, 9).Value = $blank
$wsDeDe.Cells.Item(4, 10).Value = $blank
$wsDeDe.Cells.Item(4, 11).NumberFormat = $dateFmt
$wsDeDe.Cells.Item(4, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(4, 12).Value = $blank
$wsDeDe.Cells.Item(4, 13).Value = "'True"
$wsDeDe.Cells.Item(4, 14).Value = $blank
$wsDeDe.Cells.Item(4, 15).Value = "'False"
$wsDeDe.Cells.Item(4, 16).Value = $blank

# Row 5 - e55b63a1-52a7-4daf-a112-431e7e259a2c.md
$wsDeDe.Cells.Item(5, 1).Value = "e55b63a1-52a7-4daf-a112-431e7e259a2c.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e55b63a152a74daf00000000000000000000000/e2e/e55b63a1-52a7-4daf-a112-431e7e259a2c.md", "", "", "e55b63a1-52a7-4daf-a112-431e7e259a2c.md")
$wsDeDe.Cells.Item(5, 2).Value = ".md"
$wsDeDe.Cells.Item(5, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(5, 4).Value = "e2e"
$wsDeDe.Cells.Item(5, 5).Value = "ht"
$wsDeDe.Cells.Item(5, 6).Value = "'False"
$wsDeDe.Cells.Item(5, 7).Value = "e55b63a1-52a7-4daf-a112-431e7e259a2c.f3fd377ffe290dc282c6d067ebb59dc8395379ee.de-de.xlf"
$wsDeDe.Cells.Item(5, 8).NumberFormat = $dateFmt
$wsDeDe.Cells.Item(5, 8).Value = "2016-08-12 20:49:09"
$wsDeDe.Cells.Item(5, 9).Value = $blank
$wsDeDe.Cells.Item(5, 10).Value = $blank
$wsDeDe.Cells.Item(5, 11).NumberFormat = $dateFmt
$wsDeDe.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(5, 12).Value = $blank
$wsDeDe.Cells.Item(5, 13).Value = "'True"
$wsDeDe.Cells.Item(5, 14).Value = $blank
$wsDeDe.Cells.Item(5, 15).Value = "'False"
$wsDeDe.Cells.Item(5, 16).Value = $blank
